{"js": "// Locate the target paragraphs by their stable \"Sprint N\" prefixes, then\n// perform scoped text replacements so formatting (bold) on the existing\n// runs is preserved and we don't accidentally match similar text in a\n// neighboring paragraph that happens to share the same words.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet pSprint3 = null, pSprint5 = null, pSprint6 = null, pSprint7 = null, pSprint8 = null;\nfor (const p of paragraphs.items) {\n  const t = p.text;\n  if (t.indexOf(\"Sprint 3:\") === 0) {\n    pSprint3 = p;\n  } else if (t.indexOf(\"Sprint 5\") === 0) {\n    pSprint5 = p;\n  } else if (t.indexOf(\"Sprint 6\") === 0) {\n    pSprint6 = p;\n  } else if (t.indexOf(\"Sprint 7\") === 0) {\n    pSprint7 = p;\n  } else if (t.indexOf(\"Sprint 8\") === 0) {\n    pSprint8 = p;\n  }\n}\nif (!pSprint3 || !pSprint5 || !pSprint6 || !pSprint7 || !pSprint8) {\n  throw new Error(\"Could not locate one or more target 'Sprint N' paragraphs.\");\n}\n\nasync function replaceFirst(scopeRange, searchText, replacement) {\n  const found = scopeRange.search(searchText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n  if (found.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  found.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Sprint 3: \"w.17: Window Handler and Board Menu.\" -> \"w.17: Message Service.\"\nawait replaceFirst(pSprint3, \"Window Handler and Board Menu.\", \"Message Service.\");\n\n// 2) Sprint 5: \"w.19: Message Service.\" -> \"w.19: Window Handler and Board Menu.\"\nawait replaceFirst(pSprint5, \"Message Service.\", \"Window Handler and Board Menu.\");\n\n// 3) Sprint 7: \"w.21: Contact form.\" -> \"w.21: Navigation and design.\"\nawait replaceFirst(pSprint7, \"Contact form\", \"Navigation and design\");\n\n// 4) Sprint 8: \"... (Last Iteration): Memory Game.\" -> \"... (Last Iteration): Navigation and design.\"\nawait replaceFirst(pSprint8, \"Memory Game\", \"Navigation and design\");\n\n// 5) Move the \"_GoBack\" bookmark: it used to sit inside Sprint 6 (right\n// after \"Application Nav\"); the edit removes it from there and places a\n// fresh one at the point of the most-recent edit, right before the final\n// \".\" in the Sprint 8 paragraph (i.e. immediately after \"Navigation and\n// design\").\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n} catch (e) {\n  // no-op if it doesn't exist / isn't supported\n}\n\nconst finalResults = pSprint8.search(\"Navigation and design\", { matchCase: true });\nfinalResults.load(\"items\");\nawait context.sync();\nif (finalResults.items.length === 0) {\n  throw new Error(\"Could not relocate '_GoBack' bookmark anchor text.\");\n}\nconst insertionPoint = finalResults.items[0].getRange(Word.RangeLocation.end);\ninsertionPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Apply the \"Comprehensive Project Plan\" sprint re-shuffle:\n#   Sprint 3 : \"Window Handler and Board Menu\" -> \"Message Service\"\n#   Sprint 5 : \"Message Service\"               -> \"Window Handler and Board Menu\"\n#   Sprint 7 : \"Contact form\"                  -> \"Navigation and design\"\n#   Sprint 8 : \"Memory Game\"                   -> \"Navigation and design\"\n# and move the internal \"_GoBack\" bookmark from Sprint 6 to sit right after\n# the newly-typed text in Sprint 8 (the position of the most recent edit).\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphByPrefix($doc, $prefix) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text.StartsWith($prefix)) {\n            return $p\n        }\n    }\n    throw \"Could not find a paragraph starting with '$prefix'\"\n}\n\nfunction Replace-InParagraph($paragraph, $searchText, $replacement) {\n    $rng = $paragraph.Range.Duplicate\n    $found = $rng.Find.Execute($searchText)\n    if (-not $found) {\n        throw \"Could not find text '$searchText' to replace\"\n    }\n    $rng.Text = $replacement\n    return $rng\n}\n\n# 1) Sprint 3: \"w.17: Window Handler and Board Menu.\" -> \"w.17: Message Service.\"\n$p3 = Get-ParagraphByPrefix $d \"Sprint 3:\"\nReplace-InParagraph $p3 \"Window Handler and Board Menu.\" \"Message Service.\" | Out-Null\n\n# 2) Sprint 5: \"w.19: Message Service.\" -> \"w.19: Window Handler and Board Menu.\"\n$p5 = Get-ParagraphByPrefix $d \"Sprint 5\"\nReplace-InParagraph $p5 \"Message Service.\" \"Window Handler and Board Menu.\" | Out-Null\n\n# 3) Sprint 7: \"w.21: Contact form.\" -> \"w.21: Navigation and design.\"\n$p7 = Get-ParagraphByPrefix $d \"Sprint 7\"\nReplace-InParagraph $p7 \"Contact form\" \"Navigation and design\" | Out-Null\n\n# 4) Sprint 8: \"... (Last Iteration): Memory Game.\" -> \"... (Last Iteration): Navigation and design.\"\n$p8 = Get-ParagraphByPrefix $d \"Sprint 8\"\nReplace-InParagraph $p8 \"Memory Game\" \"Navigation and design\" | Out-Null\n\n# 5) Move the \"_GoBack\" bookmark from Sprint 6 to Sprint 8 (right before the\n# trailing period, i.e. immediately after the text we just inserted).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$p8 = Get-ParagraphByPrefix $d \"Sprint 8\"\n$rng = $p8.Range.Duplicate\n$found = $rng.Find.Execute(\"Navigation and design\")\nif (-not $found) {\n    throw \"Could not relocate '_GoBack' bookmark anchor text\"\n}\n$insertionPoint = $d.Range($rng.End, $rng.End)\n$d.Bookmarks.Add(\"_GoBack\", $insertionPoint) | Out-Null\n"}
